$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDS pricing")

# --- Top parameter cells ---
$ws.Range("B1").Value = 302
$ws.Range("B2").Value = 10000000
$ws.Range("B3").Value = 0.25
$ws.Range("B4").Value = 0.05

# --- Flat 1% coupon rate instead of referencing Calibration sheet ---
$ws.Range("C6:C14").Formula = "=0.01"

# --- Insert 12 new rows (old row 15 empty spacer + totals shift down to 28-30) ---
$ws.Range("A15:A26").EntireRow.Insert()

# --- Column A: quarter offsets, continuing the +3 pattern ---
$ws.Range("A15").Formula = "=A14+3"
for ($r = 16; $r -le 25; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 1).Formula = "=A" + $prev + "+3"
}
$ws.Range("A26").Formula = "=A25+3"

# --- Columns B:L: replicate the row-14 formula pattern down through row 26 ---
for ($r = 15; $r -le 26; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 2).Formula = "=1/(1+`$B`$4/4)^(A" + $r + "/3)"
    $ws.Cells.Item($r, 3).Formula = "=0.01"
    $ws.Cells.Item($r, 4).Formula = "=D" + $prev + "*(1-C" + $prev + ")"
    $ws.Cells.Item($r, 5).Formula = "=`$B`$1/4"
    $ws.Cells.Item($r, 6).Formula = "=E" + $r + "*D" + $r + "/100"
    $ws.Cells.Item($r, 7).Formula = "=F" + $r + "*B" + $r + "*N*0.0001"
    $ws.Cells.Item($r, 8).Formula = "=D" + $prev + "*C" + $prev
    $ws.Cells.Item($r, 9).Formula = "=E" + $r + "/2*H" + $r + "/100"
    $ws.Cells.Item($r, 10).Formula = "=I" + $r + "*B" + $r + "*N*0.0001"
    $ws.Cells.Item($r, 11).Formula = "=(1-`$B`$3)*H" + $r + "/100"
    $ws.Cells.Item($r, 12).Formula = "=K" + $r + "*B" + $r + "*N"
}

# --- Totals block, now at rows 28-30, summed over the extended range ---
$ws.Range("H28").Formula = "=SUM(G7:G26)+SUM(J7:J26)"
$ws.Range("H29").Formula = "=SUM(L7:L26)"
$ws.Range("H30").Formula = "=H29-H28"

# --- Fix up the Solver defined name that pointed at the old "Value" cell ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "CDS pricing!solver_opt") {
        $n.RefersTo = "='CDS pricing'!`$H`$30"
    }
}

# --- Make "CDS pricing" the active sheet/selection ---
$ws.Activate()
$ws.Range("B2").Select()
